$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.4
$ws.Range("G2").Value = 2.44
$ws.Range("I2").Value = 3.65
$ws.Range("K2").Value = 3.3
$ws.Range("L2").Value = 1.5
$ws.Range("N2").Value = 3.15
$ws.Range("P2").Value = 1.68
$ws.Range("U2").Value = 2
$ws.Range("W2").Value = 1.54
$ws.Range("X2").Value = 10.5
$ws.Range("AD2").Value = 15.5
$ws.Range("AE2").Value = 980
$ws.Range("AK2").Value = 80
$ws.Range("AN2").Value = 28
$ws.Range("AO2").Value = 60

# Row 3
$ws.Range("F3").Value = 1.73
$ws.Range("G3").Value = 1.78
$ws.Range("I3").Value = 7.4
$ws.Range("J3").Value = 3.65
$ws.Range("K3").Value = 3.8
$ws.Range("L3").Value = 1.55
$ws.Range("N3").Value = 2.86
$ws.Range("P3").Value = 1.61
$ws.Range("T3").Value = 2.14
$ws.Range("X3").Value = 11
$ws.Range("Y3").Value = 38
$ws.Range("Z3").Value = 60
$ws.Range("AB3").Value = 6.8
$ws.Range("AF3").Value = 9.6
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 70
$ws.Range("AN3").Value = 50

# Row 4
$ws.Range("F4").Value = 2.9
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 2.84
$ws.Range("I4").Value = 2.92
$ws.Range("J4").Value = 3.15
$ws.Range("P4").Value = 1.66
$ws.Range("Q4").Value = 2.46
$ws.Range("V4").Value = 1.52
$ws.Range("W4").Value = 1.5
$ws.Range("X4").Value = 11
$ws.Range("Y4").Value = 9.800000000000001
$ws.Range("AA4").Value = 220
$ws.Range("AB4").Value = 10
$ws.Range("AC4").Value = 7.6
$ws.Range("AD4").Value = 13.5
$ws.Range("AE4").Value = 110
$ws.Range("AF4").Value = 21
$ws.Range("AG4").Value = 14
$ws.Range("AJ4").Value = 200
$ws.Range("AK4").Value = 120
$ws.Range("AL4").Value = 200
$ws.Range("AM4").Value = 580
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000

# Row 5
$ws.Range("F5").Value = 1.65
$ws.Range("K5").Value = 4.1
$ws.Range("P5").Value = 1.72
$ws.Range("U5").Value = 1.75
$ws.Range("AB5").Value = 7.8
$ws.Range("AE5").Value = 440
$ws.Range("AI5").Value = 450
$ws.Range("AM5").Value = 580

# Row 6
$ws.Range("F6").Value = 1.82
$ws.Range("G6").Value = 1.85
$ws.Range("H6").Value = 5.7
$ws.Range("I6").Value = 5.9
$ws.Range("K6").Value = 3.6
$ws.Range("P6").Value = 1.73
$ws.Range("Q6").Value = 2.1
$ws.Range("T6").Value = 1.96
$ws.Range("U6").Value = 1.83
$ws.Range("W6").Value = 2.16
$ws.Range("X6").Value = 13
$ws.Range("Y6").Value = 19
$ws.Range("Z6").Value = 1000
$ws.Range("AA6").Value = 150
$ws.Range("AB6").Value = 9.199999999999999
$ws.Range("AD6").Value = 980
$ws.Range("AI6").Value = 120
$ws.Range("AL6").Value = 48
$ws.Range("AM6").Value = 160
$ws.Range("AN6").Value = 1000
$ws.Range("AO6").Value = 120

# Row 7
$ws.Range("F7").Value = 1.09
$ws.Range("G7").Value = 1.11
$ws.Range("H7").Value = 32
$ws.Range("I7").Value = 870
$ws.Range("J7").Value = 14.5
$ws.Range("N7").Value = 7
$ws.Range("O7").Value = 1.11
$ws.Range("P7").Value = 3.4
$ws.Range("Q7").Value = 1.33
$ws.Range("R7").Value = 1.99
$ws.Range("S7").Value = 1.8
$ws.Range("U7").Value = 1.39
$ws.Range("W7").Value = 10.5
$ws.Range("Y7").Value = 980
$ws.Range("AB7").Value = 1000
$ws.Range("AD7").Value = 980
$ws.Range("AF7").Value = 9.6
$ws.Range("AJ7").Value = 9
$ws.Range("AN7").Value = 2.64

# Row 8
$ws.Range("F8").Value = 1.66
$ws.Range("G8").Value = 1.75
$ws.Range("H8").Value = 5.8
$ws.Range("I8").Value = 7.4
$ws.Range("J8").Value = 3.85
$ws.Range("K8").Value = 3.9
$ws.Range("P8").Value = 1.64
$ws.Range("Q8").Value = 2.32
$ws.Range("T8").Value = 2.28
$ws.Range("U8").Value = 1.66
$ws.Range("V8").Value = 1.16
$ws.Range("W8").Value = 2.32
$ws.Range("X8").Value = 11.5
$ws.Range("Y8").Value = 17
$ws.Range("Z8").Value = 55
$ws.Range("AA8").Value = 260
$ws.Range("AB8").Value = 6.6
$ws.Range("AC8").Value = 9.4
$ws.Range("AD8").Value = 28
$ws.Range("AE8").Value = 140
$ws.Range("AF8").Value = 8.6
$ws.Range("AG8").Value = 11
$ws.Range("AH8").Value = 32
$ws.Range("AI8").Value = 160
$ws.Range("AJ8").Value = 17
$ws.Range("AK8").Value = 22
$ws.Range("AL8").Value = 60
$ws.Range("AN8").Value = 15.5
$ws.Range("AO8").Value = 240

# Row 9
$ws.Range("G9").Value = 1.53
$ws.Range("M9").Value = 1.08
$ws.Range("Q9").Value = 1.88
$ws.Range("W9").Value = 2.88
$ws.Range("X9").Value = 1000
$ws.Range("Y9").Value = 1000
$ws.Range("AA9").Value = 1000
$ws.Range("AB9").Value = 1000
$ws.Range("AC9").Value = 1000
$ws.Range("AD9").Value = 1000
$ws.Range("AE9").Value = 1000
$ws.Range("AF9").Value = 1000
$ws.Range("AG9").Value = 1000
$ws.Range("AH9").Value = 1000
$ws.Range("AI9").Value = 1000
$ws.Range("AJ9").Value = 1000
$ws.Range("AK9").Value = 1000
$ws.Range("AL9").Value = 1000
$ws.Range("AM9").Value = 1000
$ws.Range("AN9").Value = 1000

# Row 10
$ws.Range("F10").Value = 1.38
$ws.Range("G10").Value = 1.45
$ws.Range("H10").Value = 10.5
$ws.Range("I10").Value = 14
$ws.Range("J10").Value = 4.2
$ws.Range("K10").Value = 5.4
$ws.Range("M10").Value = 1.06
$ws.Range("N10").Value = 3.9
$ws.Range("O10").Value = 1.27
$ws.Range("P10").Value = 1.99
$ws.Range("Q10").Value = 1.8
$ws.Range("R10").Value = 1.39
$ws.Range("S10").Value = 3.05
$ws.Range("T10").Value = 2.16
$ws.Range("U10").Value = 1.71
$ws.Range("V10").Value = 1.07
$ws.Range("W10").Value = 3.2
$ws.Range("X10").Value = 20
$ws.Range("Y10").Value = 32
$ws.Range("Z10").Value = 130
$ws.Range("AA10").Value = 580
$ws.Range("AB10").Value = 9.4
$ws.Range("AC10").Value = 13.5
$ws.Range("AD10").Value = 50
$ws.Range("AE10").Value = 260
$ws.Range("AF10").Value = 9.199999999999999
$ws.Range("AG10").Value = 13
$ws.Range("AH10").Value = 36
$ws.Range("AI10").Value = 200
$ws.Range("AJ10").Value = 11.5
$ws.Range("AK10").Value = 17
$ws.Range("AL10").Value = 48
$ws.Range("AM10").Value = 240
$ws.Range("AN10").Value = 8

# Row 11
$ws.Range("J11").Value = 1.03
$ws.Range("O11").Value = 1.23
$ws.Range("Q11").Value = 1.23
$ws.Range("S11").Value = 1.23
$ws.Range("T11").Value = 1.03
$ws.Range("U11").Value = 1.03

# Row 12
$ws.Range("F12").Value = 2.6
$ws.Range("G12").Value = 2.98
$ws.Range("H12").Value = 2.88
$ws.Range("I12").Value = 3.3
$ws.Range("J12").Value = 2.96
$ws.Range("K12").Value = 3.4
$ws.Range("N12").Value = 2.9
$ws.Range("P12").Value = 1.61
$ws.Range("Q12").Value = 2.48
$ws.Range("U12").Value = 1.89
$ws.Range("V12").Value = 1.44
$ws.Range("W12").Value = 1.5
$ws.Range("X12").Value = 12
$ws.Range("Y12").Value = 13.5
$ws.Range("Z12").Value = 22
$ws.Range("AB12").Value = 11
$ws.Range("AC12").Value = 8.4
$ws.Range("AD12").Value = 16.5
$ws.Range("AM12").Value = 170
